$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hide gridlines on the sheet view.
$excel.ActiveWindow.DisplayGridlines = $false

# Set explicit "best fit" column widths (matching the target workbook's
# autofit results) for columns A:F.
$ws.Columns.Item(1).ColumnWidth = 9.333333333333334
$ws.Columns.Item(2).ColumnWidth = 10.833333333333334
$ws.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws.Columns.Item(4).ColumnWidth = 16.833333333333332
$ws.Columns.Item(5).ColumnWidth = 9.333333333333334
$ws.Columns.Item(6).ColumnWidth = 9.333333333333334

# Define the print area for the sheet (adds the Print_Area defined name).
$ws.PageSetup.PrintArea = "A1:F7"

# Page setup: 2cm margins, 0.8cm header/footer, A4 paper, portrait orientation.
# (Points computed directly from cm -> pt, 2cm = 56.69291338582677pt,
# 0.8cm = 22.677165354330707pt, so the exported inch-based pageMargins
# attributes land on ~0.787400/~0.315 like the original author's file.)
$ws.PageSetup.LeftMargin = 56.69291338582677
$ws.PageSetup.RightMargin = 56.69291338582677
$ws.PageSetup.TopMargin = 56.69291338582677
$ws.PageSetup.BottomMargin = 56.69291338582677
$ws.PageSetup.HeaderMargin = 22.677165354330707
$ws.PageSetup.FooterMargin = 22.677165354330707
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
